$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial value that was bumped by one day
# (from 45189 / 2023-09-20 to 45190 / 2023-09-21) for every data row (2-398).
$ws.Range("C2:C398").Value = 45190
